$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.32218766666667
$ws.Range("H2").Value = 30.966563
$ws.Range("I2").Value = 0.06849540241430999
$ws.Range("J2").Value = 0.06849540241431
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 22.865525
$ws.Range("N2").Value = 68.596575
$ws.Range("O2").Value = 0.03261667889402277
$ws.Range("P2").Value = 0.03261667889402277
$ws.Range("Q2").Value = 236.0222401468583
$ws.Range("R2").Value = 2124.200161321725
$ws.Range("S2").Value = 0.002234092546264421
$ws.Range("T2").Value = 0.002234092546264421

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.32218766666667
$ws.Range("H3").Value = 30.966563
$ws.Range("I3").Value = 0.06849540241430999
$ws.Range("J3").Value = 0.06849540241431
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 227.1285753333333
$ws.Range("N3").Value = 681.385726
$ws.Range("O3").Value = 0.3239890537962366
$ws.Range("P3").Value = 0.3239890537962366
$ws.Range("Q3").Value = 2344.463779053304
$ws.Range("R3").Value = 21100.17401147974
$ws.Range("S3").Value = 0.02219176061760475
$ws.Range("T3").Value = 0.02219176061760476

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.32218766666667
$ws.Range("H4").Value = 30.966563
$ws.Range("I4").Value = 0.06849540241430999
$ws.Range("J4").Value = 0.06849540241431
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 207.8383993333333
$ws.Range("N4").Value = 623.5151979999999
$ws.Range("O4").Value = 0.2964724550563789
$ws.Range("P4").Value = 0.2964724550563789
$ws.Range("Q4").Value = 2145.346962258275
$ws.Range("R4").Value = 19308.12266032447
$ws.Range("S4").Value = 0.02030700011384511
$ws.Range("T4").Value = 0.02030700011384511

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 10.32218766666667
$ws.Range("H5").Value = 30.966563
$ws.Range("I5").Value = 0.06849540241430999
$ws.Range("J5").Value = 0.06849540241431
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 145.4707516666666
$ws.Range("N5").Value = 436.412255
$ws.Range("O5").Value = 0.2075077128377238
$ws.Range("P5").Value = 0.2075077128377237
$ws.Range("Q5").Value = 1501.576398714396
$ws.Range("R5").Value = 13514.18758842956
$ws.Range("S5").Value = 0.01421332429489297
$ws.Range("T5").Value = 0.01421332429489297

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 10.32218766666667
$ws.Range("H6").Value = 30.966563
$ws.Range("I6").Value = 0.06849540241430999
$ws.Range("J6").Value = 0.06849540241431
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 97.73455433333334
$ws.Range("N6").Value = 293.203663
$ws.Range("O6").Value = 0.139414099415638
$ws.Range("P6").Value = 0.1394140994156379
$ws.Range("Q6").Value = 1008.834411346697
$ws.Range("R6").Value = 9079.509702120269
$ws.Range("S6").Value = 0.009549224841702741
$ws.Range("T6").Value = 0.009549224841702741

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 95.45368733333333
$ws.Range("H7").Value = 286.361062
$ws.Range("I7").Value = 0.6334063027104161
$ws.Range("J7").Value = 0.6334063027104162
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 22.865525
$ws.Range("N7").Value = 68.596575
$ws.Range("O7").Value = 0.03261667889402277
$ws.Range("P7").Value = 0.03261667889402277
$ws.Range("Q7").Value = 2182.598674062517
$ws.Range("R7").Value = 19643.38806656265
$ws.Range("S7").Value = 0.02065960998495582
$ws.Range("T7").Value = 0.02065960998495583

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 95.45368733333333
$ws.Range("H8").Value = 286.361062
$ws.Range("I8").Value = 0.6334063027104161
$ws.Range("J8").Value = 0.6334063027104162
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 227.1285753333333
$ws.Range("N8").Value = 681.385726
$ws.Range("O8").Value = 0.3239890537962366
$ws.Range("P8").Value = 0.3239890537962366
$ws.Range("Q8").Value = 21680.26001433345
$ws.Range("R8").Value = 195122.340129001
$ws.Range("S8").Value = 0.2052167086837203
$ws.Range("T8").Value = 0.2052167086837204

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 95.45368733333333
$ws.Range("H9").Value = 286.361062
$ws.Range("I9").Value = 0.6334063027104161
$ws.Range("J9").Value = 0.6334063027104162
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 207.8383993333333
$ws.Range("N9").Value = 623.5151979999999
$ws.Range("O9").Value = 0.2964724550563789
$ws.Range("P9").Value = 0.2964724550563789
$ws.Range("Q9").Value = 19838.94158582447
$ws.Range("R9").Value = 178550.4742724203
$ws.Range("S9").Value = 0.187787521612741
$ws.Range("T9").Value = 0.187787521612741

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 95.45368733333333
$ws.Range("H10").Value = 286.361062
$ws.Range("I10").Value = 0.6334063027104161
$ws.Range("J10").Value = 0.6334063027104162
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 145.4707516666666
$ws.Range("N10").Value = 436.412255
$ws.Range("O10").Value = 0.2075077128377238
$ws.Range("P10").Value = 0.2075077128377237
$ws.Range("Q10").Value = 13885.71964573498
$ws.Range("R10").Value = 124971.4768116148
$ws.Range("S10").Value = 0.1314366931724374
$ws.Range("T10").Value = 0.1314366931724374

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 95.45368733333333
$ws.Range("H11").Value = 286.361062
$ws.Range("I11").Value = 0.6334063027104161
$ws.Range("J11").Value = 0.6334063027104162
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 97.73455433333334
$ws.Range("N11").Value = 293.203663
$ws.Range("O11").Value = 0.139414099415638
$ws.Range("P11").Value = 0.1394140994156379
$ws.Range("Q11").Value = 9329.123590996678
$ws.Range("R11").Value = 83962.1123189701
$ws.Range("S11").Value = 0.08830576925656161
$ws.Range("T11").Value = 0.08830576925656161

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 12.69470766666667
$ws.Range("H12").Value = 38.084123
$ws.Range("I12").Value = 0.08423883950185489
$ws.Range("J12").Value = 0.08423883950185491
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 22.865525
$ws.Range("N12").Value = 68.596575
$ws.Range("O12").Value = 0.03261667889402277
$ws.Range("P12").Value = 0.03261667889402277
$ws.Range("Q12").Value = 290.2711555198583
$ws.Range("R12").Value = 2612.440399678725
$ws.Range("S12").Value = 0.002747591178437122
$ws.Range("T12").Value = 0.002747591178437122

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 12.69470766666667
$ws.Range("H13").Value = 38.084123
$ws.Range("I13").Value = 0.08423883950185489
$ws.Range("J13").Value = 0.08423883950185491
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 227.1285753333333
$ws.Range("N13").Value = 681.385726
$ws.Range("O13").Value = 0.3239890537962366
$ws.Range("P13").Value = 0.3239890537962366
$ws.Range("Q13").Value = 2883.330866603144
$ws.Range("R13").Value = 25949.9777994283
$ws.Range("S13").Value = 0.027292461903099
$ws.Range("T13").Value = 0.02729246190309901

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 12.69470766666667
$ws.Range("H14").Value = 38.084123
$ws.Range("I14").Value = 0.08423883950185489
$ws.Range("J14").Value = 0.08423883950185491
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 207.8383993333333
$ws.Range("N14").Value = 623.5151979999999
$ws.Range("O14").Value = 0.2964724550563789
$ws.Range("P14").Value = 0.2964724550563789
$ws.Range("Q14").Value = 2638.447721444595
$ws.Range("R14").Value = 23746.02949300135
$ws.Range("S14").Value = 0.02497449555821519
$ws.Range("T14").Value = 0.0249744955582152

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 12.69470766666667
$ws.Range("H15").Value = 38.084123
$ws.Range("I15").Value = 0.08423883950185489
$ws.Range("J15").Value = 0.08423883950185491
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 145.4707516666666
$ws.Range("N15").Value = 436.412255
$ws.Range("O15").Value = 0.2075077128377238
$ws.Range("P15").Value = 0.2075077128377237
$ws.Range("Q15").Value = 1846.708666458596
$ws.Range("R15").Value = 16620.37799812736
$ws.Range("S15").Value = 0.017480208917134
$ws.Range("T15").Value = 0.01748020891713401

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 12.69470766666667
$ws.Range("H16").Value = 38.084123
$ws.Range("I16").Value = 0.08423883950185489
$ws.Range("J16").Value = 0.08423883950185491
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 97.73455433333334
$ws.Range("N16").Value = 293.203663
$ws.Range("O16").Value = 0.139414099415638
$ws.Range("P16").Value = 0.1394140994156379
$ws.Range("Q16").Value = 1240.711596193617
$ws.Range("R16").Value = 11166.40436574255
$ws.Range("S16").Value = 0.01174408194496957
$ws.Range("T16").Value = 0.01174408194496957

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 9.985757666666666
$ws.Range("H17").Value = 29.957273
$ws.Range("I17").Value = 0.06626293881469322
$ws.Range("J17").Value = 0.06626293881469324
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 22.865525
$ws.Range("N17").Value = 68.596575
$ws.Range("O17").Value = 0.03261667889402277
$ws.Range("P17").Value = 0.03261667889402277
$ws.Range("Q17").Value = 228.3295915711084
$ws.Range("R17").Value = 2054.966324139975
$ws.Range("S17").Value = 0.002161276997893126
$ws.Range("T17").Value = 0.002161276997893127

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 9.985757666666666
$ws.Range("H18").Value = 29.957273
$ws.Range("I18").Value = 0.06626293881469322
$ws.Range("J18").Value = 0.06626293881469324
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 227.1285753333333
$ws.Range("N18").Value = 681.385726
$ws.Range("O18").Value = 0.3239890537962366
$ws.Range("P18").Value = 0.3239890537962366
$ws.Range("Q18").Value = 2268.050912453911
$ws.Range("R18").Value = 20412.4582120852
$ws.Range("S18").Value = 0.02146846684833038
$ws.Range("T18").Value = 0.02146846684833038

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 9.985757666666666
$ws.Range("H19").Value = 29.957273
$ws.Range("I19").Value = 0.06626293881469322
$ws.Range("J19").Value = 0.06626293881469324
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 207.8383993333333
$ws.Range("N19").Value = 623.5151979999999
$ws.Range("O19").Value = 0.2964724550563789
$ws.Range("P19").Value = 0.2964724550563789
$ws.Range("Q19").Value = 2075.423889570561
$ws.Range("R19").Value = 18678.81500613505
$ws.Range("S19").Value = 0.01964513614964272
$ws.Range("T19").Value = 0.01964513614964272

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 9.985757666666666
$ws.Range("H20").Value = 29.957273
$ws.Range("I20").Value = 0.06626293881469322
$ws.Range("J20").Value = 0.06626293881469324
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 145.4707516666666
$ws.Range("N20").Value = 436.412255
$ws.Range("O20").Value = 0.2075077128377238
$ws.Range("P20").Value = 0.2075077128377237
$ws.Range("Q20").Value = 1452.635673731179
$ws.Range("R20").Value = 13073.72106358061
$ws.Range("S20").Value = 0.01375007087934302
$ws.Range("T20").Value = 0.01375007087934302

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 9.985757666666666
$ws.Range("H21").Value = 29.957273
$ws.Range("I21").Value = 0.06626293881469322
$ws.Range("J21").Value = 0.06626293881469324
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 97.73455433333334
$ws.Range("N21").Value = 293.203663
$ws.Range("O21").Value = 0.139414099415638
$ws.Range("P21").Value = 0.1394140994156379
$ws.Range("Q21").Value = 975.9535752323332
$ws.Range("R21").Value = 8783.582177090999
$ws.Range("S21").Value = 0.009237987939483976
$ws.Range("T21").Value = 0.009237987939483976

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 22.24264533333333
$ws.Range("H22").Value = 66.727936
$ws.Range("I22").Value = 0.1475965165587257
$ws.Range("J22").Value = 0.1475965165587257
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 22.865525
$ws.Range("N22").Value = 68.596575
$ws.Range("O22").Value = 0.03261667889402277
$ws.Range("P22").Value = 0.03261667889402277
$ws.Range("Q22").Value = 508.5897629354667
$ws.Range("R22").Value = 4577.3078664192
$ws.Range("S22").Value = 0.004814108186472269
$ws.Range("T22").Value = 0.004814108186472271

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 22.24264533333333
$ws.Range("H23").Value = 66.727936
$ws.Range("I23").Value = 0.1475965165587257
$ws.Range("J23").Value = 0.1475965165587257
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 227.1285753333333
$ws.Range("N23").Value = 681.385726
$ws.Range("O23").Value = 0.3239890537962366
$ws.Range("P23").Value = 0.3239890537962366
$ws.Range("Q23").Value = 5051.940346204615
$ws.Range("R23").Value = 45467.46311584154
$ws.Range("S23").Value = 0.0478196557434821
$ws.Range("T23").Value = 0.04781965574348211

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 22.24264533333333
$ws.Range("H24").Value = 66.727936
$ws.Range("I24").Value = 0.1475965165587257
$ws.Range("J24").Value = 0.1475965165587257
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 207.8383993333333
$ws.Range("N24").Value = 623.5151979999999
$ws.Range("O24").Value = 0.2964724550563789
$ws.Range("P24").Value = 0.2964724550563789
$ws.Range("Q24").Value = 4622.875803019036
$ws.Range("R24").Value = 41605.88222717132
$ws.Range("S24").Value = 0.04375830162193488
$ws.Range("T24").Value = 0.04375830162193489

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 22.24264533333333
$ws.Range("H25").Value = 66.727936
$ws.Range("I25").Value = 0.1475965165587257
$ws.Range("J25").Value = 0.1475965165587257
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 145.4707516666666
$ws.Range("N25").Value = 436.412255
$ws.Range("O25").Value = 0.2075077128377238
$ws.Range("P25").Value = 0.2075077128377237
$ws.Range("Q25").Value = 3235.654335695075
$ws.Range("R25").Value = 29120.88902125568
$ws.Range("S25").Value = 0.03062741557391638
$ws.Range("T25").Value = 0.03062741557391639

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 22.24264533333333
$ws.Range("H26").Value = 66.727936
$ws.Range("I26").Value = 0.1475965165587257
$ws.Range("J26").Value = 0.1475965165587257
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 97.73455433333334
$ws.Range("N26").Value = 293.203663
$ws.Range("O26").Value = 0.139414099415638
$ws.Range("P26").Value = 0.1394140994156379
$ws.Range("Q26").Value = 2173.875028847729
$ws.Range("R26").Value = 19564.87525962957
$ws.Range("S26").Value = 0.02057703543292003
$ws.Range("T26").Value = 0.02057703543292004
